# Auto-generated Excel COM-interop script
# Applies numeric corrections to H:N (price/profit) columns across ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H8").Value = 274.75
$ws.Range("J8").Value = 0
$ws.Range("L8").Value = 0
$ws.Range("N8").ClearContents()

$ws.Range("H17").Value = 2766075.5
$ws.Range("I17").Value = 163
$ws.Range("J17").Value = 2814600.5
$ws.Range("K17").Value = 489
$ws.Range("L17").Value = 8443801.5
$ws.Range("M17").Value = -321
$ws.Range("N17").Value = -8444137.5

$ws.Range("H113").Value = 7424
$ws.Range("I113").Value = 7424
$ws.Range("K113").Value = 7424
$ws.Range("M113").Value = -4170

$ws.Range("H132").Value = 1555.641
$ws.Range("I132").Value = 1435.1111
$ws.Range("K132").Value = 4305.3333
$ws.Range("M132").Value = -1775.3333

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H43").Value = 0
$ws.Range("J43").Value = 0
$ws.Range("L43").Value = 0
$ws.Range("N43").ClearContents()

$ws.Range("H61").Value = 6513.8335
$ws.Range("I61").Value = 7962.8887
$ws.Range("J61").Value = 2166.6667
$ws.Range("K61").Value = 7962.8887
$ws.Range("L61").Value = 2166.6667
$ws.Range("M61").Value = -7750.8887
$ws.Range("N61").Value = -2590.6667

$ws.Range("H74").Value = 3102.1636
$ws.Range("I74").Value = 3954.3635
$ws.Range("J74").Value = 1823.8636
$ws.Range("K74").Value = 3954.3635
$ws.Range("L74").Value = 1823.8636
$ws.Range("M74").Value = -3080.3635
$ws.Range("N74").Value = -3571.8636

$ws.Range("H77").Value = 3102.1636
$ws.Range("I77").Value = 3954.3635
$ws.Range("J77").Value = 1823.8636
$ws.Range("K77").Value = 19771.8175
$ws.Range("L77").Value = 9119.317999999999
$ws.Range("M77").Value = -15403.8175
$ws.Range("N77").Value = -17855.318

$ws.Range("H122").Value = 1222421
$ws.Range("I122").Value = 1509684.9
$ws.Range("J122").Value = 1549.75
$ws.Range("K122").Value = 4529054.699999999
$ws.Range("L122").Value = 4649.25
$ws.Range("M122").Value = -4526604.699999999
$ws.Range("N122").Value = -9549.25

$ws.Range("H123").Value = 37464.332
$ws.Range("J123").Value = 37464.332
$ws.Range("L123").Value = 37464.332
$ws.Range("N123").Value = -47264.332

$ws.Range("H132").Value = 5298.976
$ws.Range("I132").Value = 1718.4
$ws.Range("J132").Value = 14250.417
$ws.Range("K132").Value = 5155.200000000001
$ws.Range("L132").Value = 42751.251
$ws.Range("M132").Value = -2625.200000000001
$ws.Range("N132").Value = -47811.251

$ws.Range("H136").Value = 6513.8335
$ws.Range("I136").Value = 7962.8887
$ws.Range("J136").Value = 2166.6667
$ws.Range("K136").Value = 23888.6661
$ws.Range("L136").Value = 6500.000100000001
$ws.Range("M136").Value = -21338.6661
$ws.Range("N136").Value = -11600.0001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 3450.4033
$ws.Range("I134").Value = 3700.2444
$ws.Range("J134").Value = 2789.0588
$ws.Range("K134").Value = 11100.7332
$ws.Range("L134").Value = 8367.1764
$ws.Range("M134").Value = -8565.733200000001
$ws.Range("N134").Value = -13437.1764

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H12").Value = 2005
$ws.Range("I12").Value = 2005
$ws.Range("J12").Value = 0
$ws.Range("K12").Value = 2005
$ws.Range("L12").Value = 0
$ws.Range("M12").Value = -1835
$ws.Range("N12").ClearContents()

$ws.Range("H28").Value = 30643
$ws.Range("J28").Value = 30643
$ws.Range("L28").Value = 30643
$ws.Range("N28").Value = -31133

$ws.Range("H31").Value = 9097.538
$ws.Range("I31").Value = 4078
$ws.Range("J31").Value = 10010.182
$ws.Range("K31").Value = 4078
$ws.Range("L31").Value = 10010.182
$ws.Range("M31").Value = -3783
$ws.Range("N31").Value = -10600.182

$ws.Range("H34").Value = 9097.538
$ws.Range("I34").Value = 4078
$ws.Range("J34").Value = 10010.182
$ws.Range("K34").Value = 4078
$ws.Range("L34").Value = 10010.182
$ws.Range("M34").Value = -3876
$ws.Range("N34").Value = -10414.182

$ws.Range("H58").Value = 1639.862
$ws.Range("I58").Value = 1244.4
$ws.Range("J58").Value = 2063.5715
$ws.Range("K58").Value = 1244.4
$ws.Range("L58").Value = 2063.5715
$ws.Range("M58").Value = -1041.4
$ws.Range("N58").Value = -2469.5715

$ws.Range("H134").Value = 1281.8049
$ws.Range("I134").Value = 1234.027
$ws.Range("J134").Value = 1723.75
$ws.Range("K134").Value = 3702.081
$ws.Range("L134").Value = 5171.25
$ws.Range("M134").Value = -1167.081
$ws.Range("N134").Value = -10241.25

$ws.Range("H136").Value = 1639.862
$ws.Range("I136").Value = 1244.4
$ws.Range("J136").Value = 2063.5715
$ws.Range("K136").Value = 3733.2
$ws.Range("L136").Value = 6190.7145
$ws.Range("M136").Value = -1183.2
$ws.Range("N136").Value = -11290.7145

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H92").Value = 744.1429000000001
$ws.Range("J92").Value = 739.75
$ws.Range("L92").Value = 2219.25
$ws.Range("N92").Value = -4715.25

$ws.Range("H122").Value = 782.375
$ws.Range("J122").Value = 792.5
$ws.Range("L122").Value = 7132.5
$ws.Range("N122").Value = -12032.5

$ws.Range("H132").Value = 2158.9546
$ws.Range("J132").Value = 2233.3276
$ws.Range("L132").Value = 20099.9484
$ws.Range("N132").Value = -25159.9484

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H57").Value = 10925.444
$ws.Range("J57").Value = 10925.444
$ws.Range("L57").Value = 10925.444
$ws.Range("N57").Value = -12565.444

$ws.Range("H92").Value = 7700.25
$ws.Range("J92").Value = 7700.25
$ws.Range("L92").Value = 7700.25
$ws.Range("N92").Value = -11444.25

$ws.Range("H122").Value = 4323606.5
$ws.Range("I122").Value = 4987623
$ws.Range("K122").Value = 14962869
$ws.Range("M122").Value = -14960419

$ws.Range("H132").Value = 4110.0967
$ws.Range("I132").Value = 4713.75
$ws.Range("J132").Value = 3728.842
$ws.Range("K132").Value = 14141.25
$ws.Range("L132").Value = 11186.526
$ws.Range("M132").Value = -11611.25
$ws.Range("N132").Value = -16246.526

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H100").Value = 2026.7858
$ws.Range("I100").Value = 1527.8572
$ws.Range("J100").Value = 2525.7144
$ws.Range("K100").Value = 1527.8572
$ws.Range("L100").Value = 2525.7144
$ws.Range("M100").Value = -986.8571999999999
$ws.Range("N100").Value = -3607.7144

$ws.Range("H122").Value = 9050230
$ws.Range("J122").Value = 3335000
$ws.Range("L122").Value = 10005000
$ws.Range("N122").Value = -10009900

$ws.Range("H132").Value = 36122530
$ws.Range("I132").Value = 54181044
$ws.Range("J132").Value = 5497
$ws.Range("K132").Value = 162543132
$ws.Range("L132").Value = 16491
$ws.Range("M132").Value = -162540602
$ws.Range("N132").Value = -21551

$ws.Range("H136").Value = 7159.0977
$ws.Range("I136").Value = 5507.2256
$ws.Range("J136").Value = 12279.9
$ws.Range("K136").Value = 16521.6768
$ws.Range("L136").Value = 36839.7
$ws.Range("M136").Value = -13971.6768
$ws.Range("N136").Value = -41939.7

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H75").Value = 37971.43
$ws.Range("J75").Value = 37971.43
$ws.Range("L75").Value = 37971.43
$ws.Range("N75").Value = -39843.43

$ws.Range("H78").Value = 37971.43
$ws.Range("J78").Value = 37971.43
$ws.Range("L78").Value = 113914.29
$ws.Range("N78").Value = -123274.29

$ws.Range("H123").Value = 29436.46
$ws.Range("J123").Value = 29436.46
$ws.Range("L123").Value = 29436.46
$ws.Range("N123").Value = -39236.46

$ws.Range("H136").Value = 1703.375
$ws.Range("I136").Value = 1072.6
$ws.Range("J136").Value = 2153.9285
$ws.Range("K136").Value = 3217.8
$ws.Range("L136").Value = 6461.7855
$ws.Range("M136").Value = -667.7999999999997
$ws.Range("N136").Value = -11561.7855
